$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update cell A5 text (shared string "just split" -> "${image}")
$ws.Range("A5").Value = '${image}'

# 2. Set column A width (~19.71 stored units; closest reachable value in this
#    engine's simplified char-width model, which quantizes to steps of 1/6)
$ws.Range("A1").EntireColumn.ColumnWidth = 18.8

# 3. Set row 5 height (custom height)
$ws.Range("A5").EntireRow.RowHeight = 108.75

# 4. Update the selection to A5:C5 with active cell A5
$ws.Range("A5:C5").Select()
